{"js": "// The seller in this offer-contract template is being renamed from\n// \"\u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\" to \"\u0420\u044b\u0431\u043d\u043e\u0432 \u0420\u044b\u0431\u0430 \u0420\u044b\u0431\u043e\u0432\u0438\u0447\", and the stray\n// \"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \" phrase is removed\n// from the introductory paragraph.\n\nconst body = context.document.body;\n\n// 1) Replace every occurrence of the old full name with the new one.\n//    This covers both mentions in the intro paragraph (\"\u0418\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u044b\u0439\n//    \u043f\u0440\u0435\u0434\u043f\u0440\u0438\u043d\u0438\u043c\u0430\u0442\u0435\u043b\u044c ...\" and \"(\u0418\u041f ...)\") as well as the one in the\n//    \"6. \u0420\u0435\u043a\u0432\u0438\u0437\u0438\u0442\u044b \u043f\u0440\u043e\u0434\u0430\u0432\u0446\u0430\" section (\"\u0418\u041f \u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\").\nconst nameMatches = body.search(\"\u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\", { matchCase: true });\nnameMatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nameMatches.items.length; i++) {\n  nameMatches.items[i].insertText(\"\u0420\u044b\u0431\u043d\u043e\u0432 \u0420\u044b\u0431\u0430 \u0420\u044b\u0431\u043e\u0432\u0438\u0447\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Drop the now-obsolete \"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \"\n//    clause from the introductory paragraph.\nconst domainMatches = body.search(\"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \", { matchCase: true });\ndomainMatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < domainMatches.items.length; i++) {\n  domainMatches.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The seller in this offer-contract template is being renamed from\n# \"\u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\" to \"\u0420\u044b\u0431\u043d\u043e\u0432 \u0420\u044b\u0431\u0430 \u0420\u044b\u0431\u043e\u0432\u0438\u0447\", and the stray\n# \"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \" phrase is removed\n# from the introductory paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Replace every occurrence of the old full name with the new one.\n#    This covers both mentions in the intro paragraph (\"\u0418\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u044b\u0439\n#    \u043f\u0440\u0435\u0434\u043f\u0440\u0438\u043d\u0438\u043c\u0430\u0442\u0435\u043b\u044c ...\" and \"(\u0418\u041f ...)\") as well as the one in the\n#    \"6. \u0420\u0435\u043a\u0432\u0438\u0437\u0438\u0442\u044b \u043f\u0440\u043e\u0434\u0430\u0432\u0446\u0430\" section (\"\u0418\u041f \u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\").\n$nameFind = $d.Content.Find\n$nameFind.Text = \"\u0418\u0432\u0430\u043d\u043e\u0432 \u0418\u0432\u0430\u043d \u0418\u0432\u0430\u043d\u043e\u0432\u0438\u0447\"\n$nameFind.Replacement.Text = \"\u0420\u044b\u0431\u043d\u043e\u0432 \u0420\u044b\u0431\u0430 \u0420\u044b\u0431\u043e\u0432\u0438\u0447\"\n$nameFind.Forward = $true\n$nameFind.Wrap = 1\n$nameFind.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Drop the now-obsolete \"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \"\n#    clause from the introductory paragraph.\n$domainFind = $d.Content.Find\n$domainFind.Text = \"\u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u043e\u0433\u043e \u043d\u0430 \u0434\u043e\u043c\u0435\u043d\u0435 (www.baitstore.c\u0441\u0441), \"\n$domainFind.Replacement.Text = \"\"\n$domainFind.Forward = $true\n$domainFind.Wrap = 1\n$domainFind.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
